$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" column (F) values per row to reflect repulled data / mean
# calculation, leaving all other columns untouched.
$ws.Range("F2").Value = -1
$ws.Range("F3").Value = -1
$ws.Range("F5").Value = -1
$ws.Range("F7").Value = 2
$ws.Range("F8").Value = 3
$ws.Range("F9").Value = 3
$ws.Range("F10").Value = 10
$ws.Range("F12").Value = -2
$ws.Range("F13").Value = 3
$ws.Range("F14").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("F16").Value = -4
$ws.Range("F17").Value = 1
$ws.Range("F18").Value = 5
$ws.Range("F19").Value = 2
$ws.Range("F20").Value = 4
$ws.Range("F21").Value = 4
$ws.Range("F23").Value = 1
$ws.Range("F24").Value = -1
$ws.Range("F25").Value = -2
$ws.Range("F26").Value = -6
$ws.Range("F27").Value = -2
$ws.Range("F28").Value = 3
$ws.Range("F29").Value = -3
$ws.Range("F30").Value = -1
$ws.Range("F31").Value = -9
$ws.Range("F32").Value = -5
$ws.Range("F33").Value = 1
$ws.Range("F34").Value = 7
$ws.Range("F35").Value = 3
$ws.Range("F36").Value = 6
$ws.Range("F37").Value = -2
